# Monte Carlo results update: add EENS confidence-interval columns (Q, R)
# and refresh the simulated reliability figures for LP1-LP4 and TOTAL.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns: Q1 "EENS 95% CI", R1 "EENS 99% CI" ---
$ws.Range("Q1").Value = "EENS 95% CI"
$ws.Range("R1").Value = "EENS 99% CI"

# Copy the header formatting (bold font, border, centered) from the last
# existing header cell (P1) onto the two new header cells.
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)

# --- Row 2 (LP1) updated simulation results ---
$ws.Range("F2").Value = 2.160997846649438
$ws.Range("G2").Value = 1382
$ws.Range("H2").Value = 0.9382045643919412
$ws.Range("I2").Value = 2.303333333333333
$ws.Range("J2").Value = 483.7
$ws.Range("K2").Value = 453.8095477963819
$ws.Range("L2").Value = 197.0229585223076
$ws.Range("M2").Value = 1.156133847957449

# --- Row 3 (LP2) updated simulation results ---
$ws.Range("F3").Value = 3.153339773434606
$ws.Range("G3").Value = 1382
$ws.Range("H3").Value = 1.36903318672993
$ws.Range("I3").Value = 2.303333333333333
$ws.Range("J3").Value = 483.7
$ws.Range("K3").Value = 662.2013524212673
$ws.Range("L3").Value = 287.4969692132854
$ws.Range("M3").Value = 1.687036778787514

# --- Row 4 (LP3) updated simulation results ---
$ws.Range("F4").Value = 3.890783287857701
$ws.Range("G4").Value = 1382
$ws.Range("H4").Value = 1.689196796464993
$ws.Range("I4").Value = 2.303333333333333
$ws.Range("J4").Value = 483.7
$ws.Range("K4").Value = 817.0644904501172
$ws.Range("L4").Value = 354.7313272576486
$ws.Range("M4").Value = 2.08156905900387

# --- Row 5 (LP4) updated simulation results ---
$ws.Range("F5").Value = 4.043643784930774
$ws.Range("G5").Value = 1382
$ws.Range("H5").Value = 1.755561701127688
$ws.Range("I5").Value = 2.303333333333333
$ws.Range("J5").Value = 2.303333333333333
$ws.Range("K5").Value = 4.043643784930774
$ws.Range("L5").Value = 1.755561701127688
$ws.Range("M5").Value = 2.288702382270818

# --- Row 6 (TOTAL) updated simulation results + new EENS confidence intervals ---
$ws.Range("J6").Value = 2.303333333333333
$ws.Range("K6").Value = 3.069919230511406
$ws.Range("L6").Value = 1.332815874317542
$ws.Range("M6").Value = 7.213442068019651
$ws.Range("P6").Value = 0.04260162262065212
$ws.Range("Q6").Value = "(6.611125567935338, 7.81575856810396)"
$ws.Range("R6").Value = "(6.421826096480269, 8.005058039559028)"
